$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 15: "Print N digit binary number where 1s>=0s in all its prefixes" ----
# Date column (A) - copy date-format style from A14, then set value
$ws.Cells.Item(14, 1).Copy()
$ws.Cells.Item(15, 1).PasteSpecial(-4122)
$ws.Cells.Item(15, 1).Value = 44400

$ws.Cells.Item(15, 2).Value = "Print N digit binary number where 1s>=0s in all its prefixes"
$ws.Cells.Item(15, 3).Value = "Done - D"
$ws.Cells.Item(15, 4).Value = "Medium"

# FileName column (E) needs the wrap-text style, like B13
$ws.Cells.Item(13, 2).Copy()
$ws.Cells.Item(15, 5).PasteSpecial(-4122)
$ws.Cells.Item(15, 5).Value = "PrintNDigitBinaryNumberWhereOnesGreaterThanEqualsToZeros"

$ws.Cells.Item(15, 6).Value = "Recursion"

# Row height of 29 (two-line wrapped row, matching row 9's height)
$ws.Rows.Item(15).RowHeight = 29

# ---- Row 16: "Josephus problem (Circle of death)" ----
$ws.Cells.Item(14, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122)
$ws.Cells.Item(16, 1).Value = 44400

$ws.Cells.Item(16, 2).Value = "Josephus problem (Circle of death)"
$ws.Cells.Item(16, 3).Value = "Done - D"
$ws.Cells.Item(16, 4).Value = "Medium"
$ws.Cells.Item(16, 5).Value = "JosephusProblem"
$ws.Cells.Item(16, 6).Value = "Recursion"

# Update view/selection to match the post-edit state
$ws.Range("F16").Select() | Out-Null
